# Update the "sample final tests" sheet: fill in column E (test 3 score)
# for the rows that previously had no value recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scores = @{
    6  = 300
    7  = 700
    9  = 350
    10 = 800
    11 = 800
    12 = 800
    13 = 800
    14 = 800
    17 = 800
    18 = 583.33
    19 = 600
    21 = 750
    22 = 200
    23 = 750
    24 = 600
    25 = 400
    26 = 800
    27 = 800
    28 = 800
    29 = 300
    30 = 800
    31 = 66.67
    32 = 800
    33 = 600
    34 = 800
    35 = 800
    36 = 700
    37 = 800
    38 = 800
    39 = 800
    40 = 600
    41 = 116.67
    42 = 500
    43 = 600
    44 = 300
    45 = 350
    46 = 300
    47 = 0
    48 = 800
    49 = 750
    50 = 600
    51 = 600
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 5).Value = $scores[$row]
}
